$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.873.37'
$ws.Range("D3").Value = '1.622.25'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.57'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("E6").Value = '  -2.03%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.97'
$ws.Range("E8").Value = '  -1.95%  '
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("E10").Value = '  -1.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0881'
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("D12").Value = '1.852.83'
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").Value = '1.605.99'
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.552'
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("E16").Value = '  -1.63%  '
$ws.Range("D17").Value = '27.870.75'
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '226.66'
$ws.Range("E18").Value = '  -1.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.58'
$ws.Range("E19").Value = '  -1.39%  '
$ws.Range("D20").Value = '0.0₃0712'
$ws.Range("E20").Value = '  -1.52%  '
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.31'
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.92'
$ws.Range("E23").Value = '  -2.84%  '
$ws.Range("E24").Value = '  +1.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.09'
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("E28").Value = '  -1.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.32'
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.37'
$ws.Range("E32").Value = '  -1.01%  '
$ws.Range("D33").Value = '1.418.56'
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.07'
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("E35").Value = '  +1.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.975'
$ws.Range("E36").Value = '  -1.81%  '
$ws.Range("E37").Value = '  -0.64%  '
$ws.Range("E38").Value = '  -0.92%  '
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.843'
$ws.Range("E40").Value = '  -2.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -2.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '64.99'
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.35'
$ws.Range("E44").Value = '  -2.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.78'
$ws.Range("E45").Value = '  -3.97%  '
$ws.Range("D46").Value = '1.763.81'
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.12'
$ws.Range("E47").Value = '  -3.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.34'
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("E49").Value = '  -2.21%  '
$ws.Range("E50").Value = '  -1.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0501'
$ws.Range("E51").Value = '  -0.56%  '
